$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("H2").Value = 4.3
$ws.Range("I2").Value = 4.8
$ws.Range("K2").Value = 4.5
$ws.Range("S2").Value = 2.34
$ws.Range("AE2").Value = 48

# Row 3
$ws.Range("H3").Value = 9
$ws.Range("I3").Value = 9.199999999999999
$ws.Range("K3").Value = 5.2
$ws.Range("L3").Value = 1.31
$ws.Range("P3").Value = 2.44
$ws.Range("W3").Value = 3.25
$ws.Range("Y3").Value = 34
$ws.Range("AA3").Value = 280
$ws.Range("AL3").Value = 32
$ws.Range("AN3").Value = 5.6
$ws.Range("AO3").Value = 120

# Row 4
$ws.Range("F4").Value = 2.18
$ws.Range("G4").Value = 2.2
$ws.Range("H4").Value = 3.45
$ws.Range("I4").Value = 3.5
$ws.Range("J4").Value = 3.95
$ws.Range("K4").Value = 4
$ws.Range("M4").Value = 1.05
$ws.Range("R4").Value = 1.56
$ws.Range("S4").Value = 2.7
$ws.Range("T4").Value = 1.6
$ws.Range("V4").Value = 1.4
$ws.Range("W4").Value = 1.84
$ws.Range("X4").Value = 20
$ws.Range("Z4").Value = 27
$ws.Range("AA4").Value = 60
$ws.Range("AB4").Value = 13.5
$ws.Range("AD4").Value = 14
$ws.Range("AE4").Value = 34
$ws.Range("AF4").Value = 15.5
$ws.Range("AG4").Value = 11
$ws.Range("AH4").Value = 15
$ws.Range("AJ4").Value = 27
$ws.Range("AK4").Value = 19.5
$ws.Range("AN4").Value = 11.5
$ws.Range("AO4").Value = 25

# Row 5
$ws.Range("F5").Value = 1.51
$ws.Range("G5").Value = 1.66
$ws.Range("I5").Value = 6.4
$ws.Range("J5").Value = 4.6
$ws.Range("K5").Value = 5.7
$ws.Range("L5").Value = 1.22
$ws.Range("M5").Value = 1.02
$ws.Range("N5").Value = 5.9
$ws.Range("P5").Value = 2.9
$ws.Range("Q5").Value = 1.41
$ws.Range("R5").Value = 1.78
$ws.Range("S5").Value = 2
$ws.Range("T5").Value = 1.52
$ws.Range("U5").Value = 2.44
$ws.Range("V5").Value = 1.18
$ws.Range("W5").Value = 2.5
$ws.Range("Y5").Value = 36
$ws.Range("AE5").Value = 65
$ws.Range("AG5").Value = 11.5
$ws.Range("AH5").Value = 19
$ws.Range("AI5").Value = 60
$ws.Range("AJ5").Value = 18
$ws.Range("AK5").Value = 15.5
$ws.Range("AL5").Value = 25
$ws.Range("AM5").Value = 65
$ws.Range("AN5").Value = 5.5
$ws.Range("AO5").Value = 48

# Row 6
$ws.Range("G6").Value = 2.98
$ws.Range("I6").Value = 2.62
$ws.Range("L6").Value = 1.25
$ws.Range("R6").Value = 1.57
$ws.Range("S6").Value = 2.44
$ws.Range("U6").Value = 2.52
$ws.Range("V6").Value = 1.61

# Row 7
$ws.Range("Q7").Value = 1.72
$ws.Range("S7").Value = 2.54
$ws.Range("T7").Value = 1.84
$ws.Range("U7").Value = 2
$ws.Range("V7").Value = 1.12

# Row 8
$ws.Range("F8").Value = 2.5
$ws.Range("G8").Value = 2.54
$ws.Range("H8").Value = 3.25
$ws.Range("I8").Value = 3.35
$ws.Range("L8").Value = 1.45
$ws.Range("P8").Value = 1.81
$ws.Range("Q8").Value = 2.2
$ws.Range("T8").Value = 1.88
$ws.Range("V8").Value = 1.43
$ws.Range("W8").Value = 1.65
$ws.Range("Y8").Value = 11.5
$ws.Range("AN8").Value = 25

# Row 9
$ws.Range("F9").Value = 2.8
$ws.Range("H9").Value = 2.52
$ws.Range("I9").Value = 2.54
$ws.Range("J9").Value = 3.95
$ws.Range("K9").Value = 4
$ws.Range("R9").Value = 1.73
$ws.Range("T9").Value = 1.5
$ws.Range("U9").Value = 2.88
$ws.Range("AC9").Value = 9.6
$ws.Range("AJ9").Value = 44
$ws.Range("AL9").Value = 29

# Row 10
$ws.Range("F10").Value = 9.4
$ws.Range("G10").Value = 9.800000000000001
$ws.Range("H10").Value = 1.39
$ws.Range("I10").Value = 1.4
$ws.Range("J10").Value = 5.6
$ws.Range("K10").Value = 5.7
$ws.Range("S10").Value = 2.56
$ws.Range("V10").Value = 3.5
$ws.Range("Z10").Value = 8.800000000000001
$ws.Range("AB10").Value = 32
$ws.Range("AD10").Value = 10
$ws.Range("AF10").Value = 85
$ws.Range("AJ10").Value = 610
$ws.Range("AK10").Value = 130
$ws.Range("AN10").Value = 130

# Row 11
$ws.Range("I11").Value = 23
$ws.Range("L11").Value = 1.26
$ws.Range("P11").Value = 2.9
$ws.Range("Q11").Value = 1.5
$ws.Range("R11").Value = 1.77
$ws.Range("T11").Value = 2.36
$ws.Range("U11").Value = 1.7
$ws.Range("AC11").Value = 20
$ws.Range("AG11").Value = 13.5
$ws.Range("AH11").Value = 140

# Row 12
$ws.Range("F12").Value = 1.3
$ws.Range("G12").Value = 1.31
$ws.Range("H12").Value = 11.5
$ws.Range("I12").Value = 12
$ws.Range("L12").Value = 1.21
$ws.Range("N12").Value = 8.800000000000001
$ws.Range("O12").Value = 1.11
$ws.Range("P12").Value = 3.55
$ws.Range("Q12").Value = 1.36
$ws.Range("S12").Value = 1.9
$ws.Range("T12").Value = 1.73
$ws.Range("Z12").Value = 130
$ws.Range("AA12").Value = 380
$ws.Range("AB12").Value = 15.5
$ws.Range("AC12").Value = 17
$ws.Range("AG12").Value = 11.5
$ws.Range("AI12").Value = 100
$ws.Range("AM12").Value = 95
$ws.Range("AN12").Value = 3.3

# Row 13
$ws.Range("F13").Value = 5.8
$ws.Range("G13").Value = 5.9
$ws.Range("H13").Value = 1.67
$ws.Range("J13").Value = 4.2
$ws.Range("L13").Value = 1.38
$ws.Range("P13").Value = 2.1
$ws.Range("R13").Value = 1.44
$ws.Range("S13").Value = 3.2
$ws.Range("V13").Value = 2.46
$ws.Range("AB13").Value = 20
$ws.Range("AC13").Value = 9
$ws.Range("AJ13").Value = 150
$ws.Range("AK13").Value = 75

# Row 14
$ws.Range("G14").Value = 3.2
$ws.Range("K14").Value = 3.75
$ws.Range("P14").Value = 2.28
$ws.Range("R14").Value = 1.51
$ws.Range("T14").Value = 1.63

# Row 15
$ws.Range("F15").Value = 2.4
$ws.Range("G15").Value = 2.82
$ws.Range("H15").Value = 2.48
$ws.Range("I15").Value = 2.92
$ws.Range("J15").Value = 3.65
$ws.Range("K15").Value = 4.5
$ws.Range("L15").Value = 1.27
$ws.Range("N15").Value = 5.2
$ws.Range("O15").Value = 1.18
$ws.Range("P15").Value = 2.42
$ws.Range("Q15").Value = 1.55
$ws.Range("R15").Value = 1.58
$ws.Range("T15").Value = 1.52
$ws.Range("U15").Value = 2.5
$ws.Range("V15").Value = 1.52
$ws.Range("W15").Value = 1.54
$ws.Range("Y15").Value = 18
$ws.Range("Z15").Value = 23
$ws.Range("AA15").Value = 42
$ws.Range("AC15").Value = 10.5
$ws.Range("AD15").Value = 14
$ws.Range("AE15").Value = 27
$ws.Range("AG15").Value = 13.5
$ws.Range("AH15").Value = 16
$ws.Range("AI15").Value = 34
$ws.Range("AJ15").Value = 1000
$ws.Range("AK15").Value = 27
$ws.Range("AL15").Value = 34
$ws.Range("AM15").Value = 60
$ws.Range("AN15").Value = 15.5
$ws.Range("AO15").Value = 18

# Row 16
$ws.Range("G16").Value = 2.8
$ws.Range("H16").Value = 2.96
$ws.Range("K16").Value = 3.6
$ws.Range("V16").Value = 1.46
$ws.Range("W16").Value = 1.56
$ws.Range("X16").Value = 15
$ws.Range("Y16").Value = 13.5
$ws.Range("AA16").Value = 55
$ws.Range("AB16").Value = 12.5
$ws.Range("AC16").Value = 9
$ws.Range("AF16").Value = 20
$ws.Range("AH16").Value = 20
